$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new column B values (text, using shared strings)
# Order matters for shared string table indices, matching original authoring order
$ws.Range("B5").Value = "Audit1"
$ws.Range("B7").Value = "Admin"
$ws.Range("B6").Value = "Audit2"
$ws.Range("B9").Value = "Multi"
$ws.Range("B10").Value = "Diz"
$ws.Range("B12").Value = "Egz"
$ws.Range("B14").Value = "Dest, Elek"
$ws.Range("B15").Value = "Ism"

# Set column B width (closest achievable value to the authored 14.33203125)
$ws.Columns.Item(2).ColumnWidth = 13.43

# Update selection / view
$ws.Range("H10").Select()
